$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a value that is kept as an exact text string (not auto-converted
# to a numeric type by Excel) and leave the cell style unchanged
# (matches source cells which carry no style / "Normal" style).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "265.14"
Set-TextValue "D3" "22.85"
Set-TextValue "D4" "6.198"
Set-TextValue "D5" "0.06156"
Set-TextValue "D6" "3.560"
Set-TextValue "D7" "6.704"
Set-TextValue "D9" "0.8116"
Set-TextValue "D10" "0.1592"
Set-TextValue "D11" "0.08190"
Set-TextValue "D12" "0.03377"
Set-TextValue "D13" "0.03147"
Set-TextValue "D14" "0.09223"
Set-TextValue "D15" "3.896"
Set-TextValue "D16" "0.001713"
Set-TextValue "D17" "0.04839"
Set-TextValue "D18" "0.0006260"
Set-TextValue "D19" "0.006176"
Set-TextValue "D20" "0.006269"
Set-TextValue "D21" "0.001099"
Set-TextValue "D23" "3.697"
Set-TextValue "D26" "0.1197"
Set-TextValue "D27" "0.0002681"
Set-TextValue "D40" "0.04588"
Set-TextValue "D41" "0.006959"
Set-TextValue "D42" "0.1133"
Set-TextValue "D43" "0.003399"
Set-TextValue "D44" "0.01112"
Set-TextValue "D45" "0.00006086"
Set-TextValue "D47" "0.7699"
Set-TextValue "D48" "0.2019"

Write-Output "Applied 30 price updates"
